$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q1" (before "总计").
#    NOTE: worksheet object references can become stale/re-point once the
#    sheet collection is restructured (Add/Move/Delete), so we re-fetch the
#    sheets we need by name *after* any such structural change instead of
#    caching them beforehand.
# ---------------------------------------------------------------------------
$q1_2021 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1_2021)
$newSheet.Name = "2022-Q1"

# Re-fetch sheet references now that the collection has been restructured.
$newSheet = $wb.Worksheets.Item("2022-Q1")
$total    = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 2. Populate the header row (B1:H1) - copy the cell formatting (bold,
#    bordered, centered style) used on the "总计" sheet's header so the new
#    sheet matches that style (s="2") rather than creating a brand-new style.
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$total.Range("B1:D1").Copy()
$newSheet.Range("B1:D1").PasteSpecial(-4122)
$total.Range("B1:D1").Copy()
$newSheet.Range("E1:G1").PasteSpecial(-4122)
$total.Range("B1").Copy()
$newSheet.Range("H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Populate data rows 2-3 with the 2022-Q1 holder information.
#    Column A (index) + text columns are forced to remain text (no implicit
#    numeric conversion) by temporarily formatting as Text, then clearing
#    the format again afterwards so no stray number-format style remains.
# ---------------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2
Set-TextValue $newSheet.Range("B2") "010181"
Set-TextValue $newSheet.Range("C2") "兴业优势产业混合A"
Set-TextValue $newSheet.Range("D2") "1.07"
Set-TextValue $newSheet.Range("E2") "79.94"
Set-TextValue $newSheet.Range("F2") "3.37"
Set-TextValue $newSheet.Range("G2") "0.0361"
$newSheet.Range("H2").Value = 8

# Row 3
Set-TextValue $newSheet.Range("B3") "010182"
Set-TextValue $newSheet.Range("C3") "兴业优势产业混合C"
Set-TextValue $newSheet.Range("D3") "0.39"
Set-TextValue $newSheet.Range("E3") "79.94"
Set-TextValue $newSheet.Range("F3") "3.37"
Set-TextValue $newSheet.Range("G3") "0.0131"
$newSheet.Range("H3").Value = 8

# Column A (row index numbers), styled like the "总计" sheet's column A.
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$total.Range("A2:A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Update the "总计" (summary) sheet: insert a new first data row for
#    2022-Q1, pushing the existing 2021-Q1 row down to row 3.
# ---------------------------------------------------------------------------
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4163)
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4122)

$total.Range("A3").Value = 1

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.05

# ---------------------------------------------------------------------------
# 5. Restore the originally active/selected sheet ("2021-Q1"), since adding
#    the new worksheet shifted focus onto it.
# ---------------------------------------------------------------------------
$q1_2021 = $wb.Worksheets.Item("2021-Q1")
$q1_2021.Activate()
